$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns per latest coinranking.com scrape.
# Numeric-looking price strings are prefixed with a literal apostrophe so Excel
# stores them as text (matching the sheets existing inlineStr convention)
# instead of auto-converting them to numbers.

$ws.Range("D2").Value = '25.794.50'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '1.627.41'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '''215.32'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = '''0.5109'
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '''0.2583'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '''0.06399'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").Value = '''19.37'
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("D11").Value = '''0.07785'
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '''4.255'
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '1.627.93'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").Value = '1.851.66'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").Value = '''0.5578'
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").Value = '''63.49'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '0.0₅7536'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("D18").Value = '25.798.54'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '''1.004'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '''193.91'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("E21").Value = '  -2.94%  '
$ws.Range("D22").Value = '''9.797'
$ws.Range("E22").Value = '  -1.91%  '
$ws.Range("D23").Value = '''6.001'
$ws.Range("E23").Value = '  -1.59%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = '''1.808'
$ws.Range("E25").Value = '  -5.15%  '
$ws.Range("D26").Value = '''0.1306'
$ws.Range("E26").Value = '  +5.60%  '
$ws.Range("D27").Value = '''141.45'
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").Value = '''6.741'
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("D30").Value = '''1.238'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").Value = '''0.04884'
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").Value = '''3.294'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '''3.187'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("D35").Value = '''2.377'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '''0.8955'
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("D37").Value = '1.136.02'
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("D38").Value = '''2.542'
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").Value = '''0.5489'
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").Value = '''1.002'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").Value = '''5.595'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").Value = '''0.7955'
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").Value = '''97.35'
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("D45").Value = '1.775.28'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").Value = '  -7.83%  '
$ws.Range("D47").Value = '''0.4426'
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("D48").Value = '''54.88'
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").Value = '''0.05068'
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("D50").Value = '''7.553'
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("D51").Value = '''1.006'
$ws.Range("E51").Value = '  -0.11%  '
